$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update the Runmode for the "Registration" test case (row 2) from YES to NO
$ws.Range("C2").Value = "NO"

# Update the active selection to reflect where the user continued editing
$ws.Range("B8").Select()
